$wb = $excel.ActiveWorkbook

# Sheet "展览": update F4 (1346 -> 1358) and F5 (651 -> 653)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1358
$ws1.Range("F5").Value = 653

# Sheet "全部类型": update F4 (1346 -> 1358) and F6 (651 -> 653)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1358
$ws4.Range("F6").Value = 653
